$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 850
$ws.Range("I18").Value = 836
$ws.Range("K18").Value = 836
$ws.Range("M18").Value = -552

# Row 51
$ws.Range("H51").Value = 78333
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 78333
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 78333
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -79301

# Row 112
$ws.Range("H112").Value = 4611.25
$ws.Range("J112").Value = 4782.421
$ws.Range("L112").Value = 14347.263
$ws.Range("N112").Value = -16563.263

# Row 125
$ws.Range("H125").Value = 7920
$ws.Range("I125").Value = 7920
$ws.Range("K125").Value = 71280
$ws.Range("M125").Value = -68820

# Row 132
$ws.Range("H132").Value = 1688.659
$ws.Range("I132").Value = 1719.4875
$ws.Range("J132").Value = 1380.375
$ws.Range("K132").Value = 5158.4625
$ws.Range("L132").Value = 4141.125
$ws.Range("M132").Value = -2628.4625
$ws.Range("N132").Value = -9201.125

# Row 134
$ws.Range("H134").Value = 88000
$ws.Range("J134").Value = 88000
$ws.Range("L134").Value = 88000
$ws.Range("N134").Value = -98140

# Row 137
$ws.Range("H137").Value = 2032.4849
$ws.Range("I137").Value = 1873.3334
$ws.Range("J137").Value = 2311
$ws.Range("K137").Value = 5620.0002
$ws.Range("L137").Value = 6933
$ws.Range("M137").Value = -3070.0002
$ws.Range("N137").Value = -12033

# Row 138
$ws.Range("H138").Value = 2747.5676
$ws.Range("J138").Value = 5378.793
$ws.Range("L138").Value = 16136.379
$ws.Range("N138").Value = -26416.379

$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Range("H10").Value = 10000
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -9830
$ws.Range("N10").ClearContents()

# Row 45
$ws.Range("H45").Value = 2552.4138
$ws.Range("I45").Value = 2465
$ws.Range("K45").Value = 2465
$ws.Range("M45").Value = -2088

# Row 61
$ws.Range("H61").Value = 6364.294
$ws.Range("I61").Value = 3323.077
$ws.Range("K61").Value = 3323.077
$ws.Range("M61").Value = -3111.077

# Row 102
$ws.Range("H102").Value = 2231.5
$ws.Range("I102").Value = 2156.6428
$ws.Range("K102").Value = 2156.6428
$ws.Range("M102").Value = -534.6428000000001

# Row 122
$ws.Range("H122").Value = 4126.8887
$ws.Range("I122").Value = 4085.6
$ws.Range("J122").Value = 4333.3335
$ws.Range("K122").Value = 12256.8
$ws.Range("L122").Value = 13000.0005
$ws.Range("M122").Value = -9806.799999999999
$ws.Range("N122").Value = -17900.0005

# Row 136
$ws.Range("H136").Value = 6364.294
$ws.Range("I136").Value = 3323.077
$ws.Range("K136").Value = 9969.231
$ws.Range("M136").Value = -7419.231

$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 41879.832
$ws.Range("J81").Value = 41879.832
$ws.Range("L81").Value = 41879.832
$ws.Range("N81").Value = -44001.832

# Row 84
$ws.Range("H84").Value = 41879.832
$ws.Range("J84").Value = 41879.832
$ws.Range("L84").Value = 125639.496
$ws.Range("N84").Value = -136247.496

# Row 107
$ws.Range("H107").Value = 2155.68
$ws.Range("I107").Value = 1591.2632
$ws.Range("J107").Value = 3943
$ws.Range("K107").Value = 1591.2632
$ws.Range("L107").Value = 3943
$ws.Range("M107").Value = 328.7367999999999
$ws.Range("N107").Value = -7783

# Row 134
$ws.Range("H134").Value = 2692.625
$ws.Range("I134").Value = 2650.484
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 7951.451999999999
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -5416.451999999999
$ws.Range("N134").Value = -17067

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6541.032
$ws.Range("I31").Value = 6228.8823
$ws.Range("K31").Value = 6228.8823
$ws.Range("M31").Value = -5933.8823

# Row 34
$ws.Range("H34").Value = 6541.032
$ws.Range("I34").Value = 6228.8823
$ws.Range("K34").Value = 6228.8823
$ws.Range("M34").Value = -6026.8823

# Row 105
$ws.Range("H105").Value = 2627.9167
$ws.Range("I105").Value = 2627.9167
$ws.Range("K105").Value = 2627.9167
$ws.Range("M105").Value = -880.9167000000002

# Row 132
$ws.Range("H132").Value = 1751
$ws.Range("I132").Value = 1751
$ws.Range("K132").Value = 5253
$ws.Range("M132").Value = -2723

# Row 134
$ws.Range("H134").Value = 1732.2927
$ws.Range("I134").Value = 1728.2632
$ws.Range("K134").Value = 5184.7896
$ws.Range("M134").Value = -2649.7896

$ws = $wb.Worksheets.Item("CUL")
# Row 60
$ws.Range("H60").Value = 1599.6666
$ws.Range("I60").Value = 2800
$ws.Range("J60").Value = 999.5
$ws.Range("K60").Value = 8400
$ws.Range("L60").Value = 2998.5
$ws.Range("M60").Value = -8149
$ws.Range("N60").Value = -3500.5

# Row 116
$ws.Range("H116").Value = 2263.6667
$ws.Range("I116").Value = 2263.6667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 6791.000100000001
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -3349.000100000001
$ws.Range("N116").ClearContents()

# Row 122
$ws.Range("H122").Value = 4656.5454
$ws.Range("I122").Value = 1593.3334
$ws.Range("J122").Value = 8332.4
$ws.Range("K122").Value = 14340.0006
$ws.Range("L122").Value = 74991.59999999999
$ws.Range("M122").Value = -11890.0006
$ws.Range("N122").Value = -79891.59999999999

# Row 136
$ws.Range("H136").Value = 5631.8423
$ws.Range("I136").Value = 4952
$ws.Range("J136").Value = 9257.666999999999
$ws.Range("K136").Value = 14856
$ws.Range("L136").Value = 27773.001
$ws.Range("M136").Value = -9756
$ws.Range("N136").Value = -37973.001

$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 1402000.2
$ws.Range("I7").Value = 2005001
$ws.Range("J7").Value = 999999.7
$ws.Range("K7").Value = 2005001
$ws.Range("L7").Value = 999999.7
$ws.Range("M7").Value = -2004889
$ws.Range("N7").Value = -1000223.7

# Row 8
$ws.Range("H8").Value = 1402000.2
$ws.Range("I8").Value = 2005001
$ws.Range("J8").Value = 999999.7
$ws.Range("K8").Value = 2005001
$ws.Range("L8").Value = 999999.7
$ws.Range("M8").Value = -2004862
$ws.Range("N8").Value = -1000277.7

# Row 12
$ws.Range("H12").Value = 500000
$ws.Range("I12").Value = 500000
$ws.Range("K12").Value = 500000
$ws.Range("M12").Value = -499860

# Row 18
$ws.Range("H18").Value = 213000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 213000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 213000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -213586

# Row 20
$ws.Range("H20").Value = 17500
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 21666.666
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 21666.666
$ws.Range("M20").Value = -4755
$ws.Range("N20").Value = -22156.666

# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

# Row 48
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30970

# Row 102
$ws.Range("H102").Value = 2453.1667
$ws.Range("I102").Value = 2159.889
$ws.Range("K102").Value = 2159.889
$ws.Range("M102").Value = -537.8890000000001

# Row 132
$ws.Range("H132").Value = 4849.28
$ws.Range("I132").Value = 4362.25
$ws.Range("K132").Value = 13086.75
$ws.Range("M132").Value = -10556.75

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2124.275
$ws.Range("I40").Value = 2119.8057
$ws.Range("K40").Value = 2119.8057
$ws.Range("M40").Value = -1983.8057

# Row 122
$ws.Range("H122").Value = 4494.5454
$ws.Range("J122").Value = 9748.75
$ws.Range("L122").Value = 29246.25
$ws.Range("N122").Value = -34146.25

# Row 132
$ws.Range("H132").Value = 7291.25
$ws.Range("I132").Value = 6139.5
$ws.Range("J132").Value = 13050
$ws.Range("K132").Value = 18418.5
$ws.Range("L132").Value = 39150
$ws.Range("M132").Value = -15888.5
$ws.Range("N132").Value = -44210

# Row 136
$ws.Range("H136").Value = 1442.42
$ws.Range("I136").Value = 1462.9375
$ws.Range("J136").Value = 950
$ws.Range("K136").Value = 4388.8125
$ws.Range("L136").Value = 2850
$ws.Range("M136").Value = -1838.8125
$ws.Range("N136").Value = -7950

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# Row 9
$ws.Range("H9").Value = 4000
$ws.Range("I9").Value = 4000
$ws.Range("K9").Value = 4000
$ws.Range("M9").Value = -3860

# Row 14
$ws.Range("H14").Value = 1019.0476
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 1800
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1800
$ws.Range("M14").Value = -832
$ws.Range("N14").Value = -2136

# Row 45
$ws.Range("H45").Value = 99991.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 99991.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 99991.5
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -100973.5

# Row 132
$ws.Range("H132").Value = 4433.8076
$ws.Range("I132").Value = 4335.409
$ws.Range("J132").Value = 4975
$ws.Range("K132").Value = 13006.227
$ws.Range("L132").Value = 14925
$ws.Range("M132").Value = -10476.227
$ws.Range("N132").Value = -19985

# Row 136
$ws.Range("H136").Value = 2965.5854
$ws.Range("I136").Value = 3020.5625
$ws.Range("K136").Value = 9061.6875
$ws.Range("M136").Value = -6511.6875

